# Sample Data for Sentiment Analysis.xlsx — "Add files via upload" re-save.
#
# Observed changes (from the canonical-XML diff):
#   1. Worksheet "Sheet1" renamed to "Responses".
#      (The workbook-scoped hidden defined name _xlnm._FilterDatabase
#       tracks the sheet it points at and is updated automatically by
#       Excel when the sheet is renamed.)
#   2. The cached UI selection on the sheet (activeCell/sqref "C24") is
#      cleared back to the default top-left cell, A1.
#
# (The other lines in the diff — the x15ac:absPath "last saved from"
# folder and the xr:revisionPtr documentId GUID — are Excel-internal,
# session-specific bookkeeping values that aren't exposed on the
# Workbook/Worksheet COM object model; they're stamped by the host
# application itself on save, not by user/automation code, so there is
# nothing to set for them here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet (also re-points the _FilterDatabase defined name).
$ws.Name = "Responses"

# 2. Reset the saved selection to A1 instead of the previous C24.
$ws.Range("A1").Select()
